# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E) and "Correspond Handback DateTime" (H)
# timestamps for row 2 (the ae18194b... file) on both the zh-cn and de-de sheets,
# reflecting a fresh handback report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 06:06:00"
$wsZhCn.Range("H2").Value = "2016-03-17 06:06:40"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 06:06:08"
$wsDeDe.Range("H2").Value = "2016-03-17 06:06:53"
